$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels I0 and IF
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting (bold, border, centered) used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the I0 / IF data values for rows 2-16
$data = @(
    @(10, 10),
    @(7, 8),
    @(6, 6),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(1, 3),
    @(4, 8),
    @(3, 7),
    @(5, 6),
    @(1, 2),
    @(6, 8),
    @(1, 4),
    @(1, 4),
    @(3, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
